$wb = $excel.ActiveWorkbook

# Bump the form_version on the "settings" sheet (B3), leaving the cursor
# on B4 afterwards - mirrors pressing Enter after typing the new value.
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Activate() | Out-Null
$wsSettings.Range("B3").Value = 20210421001
$wsSettings.Range("B4").Select() | Out-Null

# Mark the geopoint question ("stand_outside") as required on the
# "survey" sheet (C2), leaving the cursor on C3 afterwards.
$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Activate() | Out-Null
$wsSurvey.Range("C2").Value = 1
$wsSurvey.Range("C3").Select() | Out-Null
